$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (rows 449-451), pushing the
# existing rows 449-537 down to 452-540.
$ws.Rows("449:451").Insert()

# Row 449 (new)
$ws.Range("A449").Value = 8
$ws.Range("B449").Value = "Terminal La Palmera de La Serena"
$ws.Range("C449").Value = "Coquimbo"
$ws.Range("D449").Value = 44637
$ws.Range("E449").Value = 4
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100104
$ws.Range("H449").Value = "Frutos de pepita"
$ws.Range("I449").Value = 100104005
$ws.Range("J449").Value = "Pera"
$ws.Range("K449").Value = "Packham's Triumph"
$ws.Range("L449").Value = "Especial"
$ws.Range("M449").Value = 16
$ws.Range("N449").Value = 245000
$ws.Range("O449").Value = 250000
$ws.Range("P449").Value = 247500
$ws.Range("Q449").Value = "`$/bins (450 kilos)"
$ws.Range("R449").Value = "Región de O'Higgins"
$ws.Range("S449").Value = 550
$ws.Range("T449").Value = 450

# Row 450 (new)
$ws.Range("A450").Value = 8
$ws.Range("B450").Value = "Terminal La Palmera de La Serena"
$ws.Range("C450").Value = "Coquimbo"
$ws.Range("D450").Value = 44637
$ws.Range("E450").Value = 4
$ws.Range("F450").Value = "Fruta"
$ws.Range("G450").Value = 100104
$ws.Range("H450").Value = "Frutos de pepita"
$ws.Range("I450").Value = 100104005
$ws.Range("J450").Value = "Pera"
$ws.Range("K450").Value = "Packham's Triumph"
$ws.Range("L450").Value = "Primera"
$ws.Range("M450").Value = 20
$ws.Range("N450").Value = 215000
$ws.Range("O450").Value = 220000
$ws.Range("P450").Value = 217500
$ws.Range("Q450").Value = "`$/bins (450 kilos)"
$ws.Range("R450").Value = "Región de O'Higgins"
$ws.Range("S450").Value = 483
$ws.Range("T450").Value = 450

# Row 451 (new)
$ws.Range("A451").Value = 8
$ws.Range("B451").Value = "Terminal La Palmera de La Serena"
$ws.Range("C451").Value = "Coquimbo"
$ws.Range("D451").Value = 44637
$ws.Range("E451").Value = 4
$ws.Range("F451").Value = "Fruta"
$ws.Range("G451").Value = 100104
$ws.Range("H451").Value = "Frutos de pepita"
$ws.Range("I451").Value = 100104005
$ws.Range("J451").Value = "Pera"
$ws.Range("K451").Value = "Packham's Triumph"
$ws.Range("L451").Value = "Segunda"
$ws.Range("M451").Value = 20
$ws.Range("N451").Value = 185000
$ws.Range("O451").Value = 190000
$ws.Range("P451").Value = 187500
$ws.Range("Q451").Value = "`$/bins (450 kilos)"
$ws.Range("R451").Value = "Región de O'Higgins"
$ws.Range("S451").Value = 417
$ws.Range("T451").Value = 450
